$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Choose Courier by",
  "Shipping Insurance",
  "Taxes & Duties Paid by*",
  "Platform",
  "Platform Order Number",
  "Order Tags",
  "Receiver's Full Name*",
  "Receiver's Phone Number*",
  "Receover's Email",
  "Receiver's Tax ID",
  "Receiver's Address Line 1*",
  "Receiver's Address Line 2",
  "Receiver's Postal Code*",
  "Receiver's City*",
  "Receiver's State/Province",
  "Receiver's Country*",
  "Item Length (in)*",
  "Item Width (in)*",
  "Item Height (in)*",
  "Item Weight (lb)*",
  "Item Category/HS codes*",
  "Item Contains Liquid",
  "Item Contains Battery",
  "Item Description*",
  "Item Country of Origin",
  "Item SKU",
  "Item Customs Value*",
  "Item Customs Value Currency*",
  "Item Quantity",
  "Buyer's Notes",
  "Seller's Notes"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Cells.Item(2, 24).Value = "LIC Football Merch"
